$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Values
$ws.Range("B1").Value = 0
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "disconnected_elements"

# Format B1: bold font, thin box border, centered horizontally, top vertically
$b1 = $ws.Range("B1")
$b1.Font.Bold = $true
$b1.HorizontalAlignment = -4108  # xlCenter
$b1.VerticalAlignment = -4160    # xlTop
$b1.Borders.LineStyle = 1        # xlContinuous
$b1.Borders.Weight = 2           # xlThin

# Copy B1's formatting onto A2 so both cells share the same style entry
# (rather than re-deriving it property-by-property, which would mint a
# redundant/duplicate style record).
$b1.Copy()
$a2 = $ws.Range("A2")
$a2.PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
